$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 865
$ws.Range("F3").Value = 1800
$ws.Range("F4").Value = 84
$ws.Range("F8").Value = 2119
$ws.Range("F10").Value = 51
$ws.Range("F12").Value = 669
$ws.Range("F14").Value = 4040
$ws.Range("F16").Value = 376
$ws.Range("F17").Value = 3154
$ws.Range("F18").Value = 861
$ws.Range("F21").Value = 184
$ws.Range("F22").Value = 2103
$ws.Range("F24").Value = 4
$ws.Range("F25").Value = 1993
$ws.Range("F26").Value = 402
$ws.Range("F27").Value = 219
$ws.Range("F29").Value = 8809
$ws.Range("F33").Value = 771
$ws.Range("F34").Value = 19
$ws.Range("F38").Value = 949
$ws.Range("F41").Value = 200
$ws.Range("F42").Value = 175
$ws.Range("F43").Value = 4676
$ws.Range("F46").Value = 88
$ws.Range("F47").Value = 418

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 90
$ws.Range("F17").Value = 3438

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 8417
$ws.Range("F4").Value = 1352

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 865
$ws.Range("F4").Value = 1352
$ws.Range("F5").Value = 84
$ws.Range("F7").Value = 1382
$ws.Range("F10").Value = 51
$ws.Range("F12").Value = 4040
$ws.Range("F14").Value = 376
$ws.Range("F15").Value = 3154
$ws.Range("F16").Value = 861
$ws.Range("F19").Value = 2103
$ws.Range("F24").Value = 4
$ws.Range("F25").Value = 1993
$ws.Range("F27").Value = 402
$ws.Range("F28").Value = 219
$ws.Range("F30").Value = 8809
$ws.Range("F34").Value = 771
$ws.Range("F36").Value = 949
$ws.Range("F39").Value = 200
$ws.Range("F41").Value = 175
$ws.Range("F42").Value = 4676
$ws.Range("F45").Value = 88
$ws.Range("F46").Value = 418
